# The edit rotates the data of rows 3-7 (sheet "Artfynd") up by one:
#   row4 -> row3, row5 -> row4, row6 -> row5, row7 -> row6, row3 -> row7
# (row-invariant columns such as C, T, U, V, W, Z, AB, AD, AE, AG are
# identical across these rows and therefore need no change).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $text) {
    # Force text storage so date-like strings ("2016-08-02") aren't
    # reinterpreted by Excel as date serial numbers.
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $text
}

# --- Row 3 (becomes old row 4's data) ---
$ws.Range("A3").Value = 65011367
$ws.Range("B3").Value = 78527
$ws.Range("D3").Value = "LC"
$ws.Range("E3").Value = 229497
$ws.Range("F3").Value = "Korallblylav"
$ws.Range("G3").Value = "Parmeliella triptophylla"
$ws.Range("H3").Value = "(Ach.) Müll.Arg."
$ws.Range("Q3").Value = 430899.7582483087
$ws.Range("R3").Value = 7050903.864986234

# --- Row 4 (becomes old row 5's data) ---
$ws.Range("A4").Value = 65011365
$ws.Range("B4").Value = 78570
$ws.Range("D4").Value = "NT"
$ws.Range("E4").Value = 2081
$ws.Range("F4").Value = "Skrovellav"
$ws.Range("G4").Value = "Lobaria scrobiculata"
$ws.Range("H4").Value = "(Scop.) DC."
$ws.Range("Q4").Value = 431021.4584440839
$ws.Range("R4").Value = 7050920.376071113

# --- Row 5 (becomes old row 6's data) ---
$ws.Range("A5").Value = 65011364
$ws.Range("B5").Value = 78569
$ws.Range("E5").Value = 6458
$ws.Range("F5").Value = "Lunglav"
$ws.Range("G5").Value = "Lobaria pulmonaria"
$ws.Range("H5").Value = "(L.) Hoffm."
$ws.Range("Q5").Value = 430989.0331701299
$ws.Range("R5").Value = 7050947.365669774
Set-TextCell "Y5" "2016-07-22"
Set-TextCell "AA5" "2016-07-22"

# --- Row 6 (becomes old row 7's data) ---
$ws.Range("A6").Value = 101328774
$ws.Range("B6").Value = 78527
$ws.Range("D6").Value = "LC"
$ws.Range("E6").Value = 229497
$ws.Range("F6").Value = "Korallblylav"
$ws.Range("G6").Value = "Parmeliella triptophylla"
$ws.Range("H6").Value = "(Ach.) Müll.Arg."
$ws.Range("J6").Value = ""
$ws.Range("K6").Value = ""
$ws.Range("N6").Value = ""
$ws.Range("P6").Value = "Getbovägen , Jmt"
$ws.Range("Q6").Value = 430901.891005214
$ws.Range("R6").Value = 7050878.869665774
$ws.Range("S6").Value = 10
Set-TextCell "Y6" "2022-05-31"
Set-TextCell "AA6" "2022-05-31"
$ws.Range("AF6").Value = ""
$ws.Range("AW6").Value = "Erik Söderhjelm"
$ws.Range("AX6").Value = "Erik Söderhjelm"
$ws.Range("AY6").ClearContents()

# --- Row 7 (becomes old row 3's data) ---
$ws.Range("A7").Value = 65011366
$ws.Range("B7").Value = 78570
$ws.Range("D7").Value = "NT"
$ws.Range("E7").Value = 2081
$ws.Range("F7").Value = "Skrovellav"
$ws.Range("G7").Value = "Lobaria scrobiculata"
$ws.Range("H7").Value = "(Scop.) DC."
$ws.Range("J7").ClearContents()
$ws.Range("K7").ClearContents()
$ws.Range("N7").ClearContents()
$ws.Range("P7").Value = "N Glasbruksberget, Jmt"
$ws.Range("Q7").Value = 431108.297595304
$ws.Range("R7").Value = 7050707.757188959
$ws.Range("S7").Value = 25
Set-TextCell "Y7" "2016-08-02"
Set-TextCell "AA7" "2016-08-02"
$ws.Range("AF7").ClearContents()
$ws.Range("AW7").Value = "Sebastian Acker"
$ws.Range("AX7").Value = "Sebastian Acker"
$ws.Range("AY7").Value = "SCA Skog Naturvärdesinventering"
